# Weekly update: a new "Primera" quality price-report row for the current
# week is inserted at the top of the data (row 5, right after the three
# most-recent pre-existing "Primera" entries), pushing every later row down
# by one. The new row reuses the same market/product/quality metadata as
# the row immediately below it (which used to be row 5) and only carries
# fresh Fecha / Volumen / Precio mínimo / Precio máximo / Precio promedio
# ponderado / Precio $/Kg figures for this week's report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5; Excel shifts rows 5..134 down to 6..135
# and grows the used range (dimension) automatically.
$ws.Rows.Item(5).Insert()

# Seed the new row with the same descriptive columns as the row right
# below it (the old row 5), so we don't need to retype any text (and avoid
# any accented-character / encoding pitfalls along the way).
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(5, $col).Value = $ws.Cells.Item(6, $col).Value2
}

# Overwrite just the figures that are new for this week's report.
$ws.Cells.Item(5, 4).Value  = 44552  # D: Fecha
$ws.Cells.Item(5, 13).Value = 80     # M: Volumen
$ws.Cells.Item(5, 14).Value = 8500   # N: Precio mínimo
$ws.Cells.Item(5, 15).Value = 9000   # O: Precio máximo
$ws.Cells.Item(5, 16).Value = 8750   # P: Precio promedio ponderado
$ws.Cells.Item(5, 19).Value = 2188   # S: Precio $/Kg
